$p = $ppt.ActivePresentation
$s = $p.Slides.Item(26)

# Move the existing title textbox ("Design of a static website")
$title = $s.Shapes.Item(1)
$title.Left = 230.79614173228347
$title.Top = 58.046535433070865

# Burn shape ids 2 and 3 (deleted scratch shapes) so the real textbox we
# add below lands on id=4, matching the authored deck.
$scratch1 = $s.Shapes.AddTextbox(1, 0, 0, 1, 1)
$scratch2 = $s.Shapes.AddTextbox(1, 0, 0, 1, 1)
$scratch1.Delete()
$scratch2.Delete()

# Add the new textbox with the course URL
$box = $s.Shapes.AddTextbox(1, 165.32047244094488, 227.58984251968505, 648.6989763779527, 117.58725)
$box.Name = "TextBox 3"

$tf = $box.TextFrame
$tf.WordWrap = $true
$tf.AutoSize = 1

$tr = $tf.TextRange
$tr.Text = "https://github.com/jianchentech/WebCamp/blob/main/WebProject/WebStatic.md"
$tr.LanguageID = "en-CA"
$tr.Font.Size = 32
$tr.Font.Color.RGB = 255

$tr2 = $box.TextFrame2.TextRange
$tr2.ParagraphFormat.SpaceWithin = 1500

# Restore the authored height (spAutoFit relayout differs from the
# runtime's text metrics) and re-assert noFill last, since touching
# AutoSize/Height drops it.
$box.Height = 117.58725
$box.Fill.Visible = $false
